# Auto-applies the cell-level numeric updates captured in the commit diff
# across the ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 21
$ws.Range("H21").Value = 8672
$ws.Range("I21").Value = 8672
$ws.Range("K21").Value = 8672
$ws.Range("M21").Value = -8204
# Row 23
$ws.Range("H23").Value = 8672
$ws.Range("I23").Value = 8672
$ws.Range("K23").Value = 8672
$ws.Range("M23").Value = -8438
# Row 86
$ws.Range("H86").Value = 4584.1
$ws.Range("I86").Value = 5736.3335
$ws.Range("J86").Value = 3815.9443
$ws.Range("K86").Value = 5736.3335
$ws.Range("L86").Value = 3815.9443
$ws.Range("M86").Value = -4613.3335
$ws.Range("N86").Value = -6061.9443
# Row 89
$ws.Range("H89").Value = 4584.1
$ws.Range("I89").Value = 5736.3335
$ws.Range("J89").Value = 3815.9443
$ws.Range("K89").Value = 28681.6675
$ws.Range("L89").Value = 19079.7215
$ws.Range("M89").Value = -23065.6675
$ws.Range("N89").Value = -30311.7215
# Row 111
$ws.Range("H111").Value = 1636.35
$ws.Range("I111").Value = 985.7273
$ws.Range("K111").Value = 2957.1819
$ws.Range("M111").Value = 109.8181
# Row 139
$ws.Range("H139").Value = 0
$ws.Range("J139").Value = 0
$ws.Range("L139").Value = 0
$ws.Range("N139").ClearContents()
# Row 140
$ws.Range("H140").Value = 97500
$ws.Range("J140").Value = 97500
$ws.Range("L140").Value = 97500
$ws.Range("N140").Value = -107860

$ws = $wb.Worksheets.Item("ARM")
# Row 131
$ws.Range("H131").Value = 0
$ws.Range("I131").Value = 0
$ws.Range("K131").Value = 0
$ws.Range("M131").ClearContents()
# Row 132
$ws.Range("H132").Value = 3313.8333
$ws.Range("I132").Value = 2546.147
$ws.Range("K132").Value = 7638.441
$ws.Range("M132").Value = -5108.441
# Row 134
$ws.Range("H134").Value = 121427.14
$ws.Range("J134").Value = 121427.14
$ws.Range("L134").Value = 121427.14
$ws.Range("N134").Value = -131567.14
# Row 135
$ws.Range("H135").Value = 181373.5
$ws.Range("J135").Value = 181373.5
$ws.Range("L135").Value = 181373.5
$ws.Range("N135").Value = -191513.5

$ws = $wb.Worksheets.Item("BSM")
# Row 107
$ws.Range("H107").Value = 2454.5293
$ws.Range("I107").Value = 2091.4814
$ws.Range("J107").Value = 3854.8572
$ws.Range("K107").Value = 2091.4814
$ws.Range("L107").Value = 3854.8572
$ws.Range("M107").Value = -171.4814000000001
$ws.Range("N107").Value = -7694.8572
# Row 132
$ws.Range("H132").Value = 138230.22
$ws.Range("J132").Value = 138230.22
$ws.Range("L132").Value = 138230.22
$ws.Range("N132").Value = -148350.22
# Row 134
$ws.Range("H134").Value = 2901.1562
$ws.Range("I134").Value = 2135.4375
$ws.Range("K134").Value = 6406.3125
$ws.Range("M134").Value = -3871.3125

$ws = $wb.Worksheets.Item("CRP")
# Row 22
$ws.Range("H22").Value = 9957.583
$ws.Range("I22").Value = 11861.1
$ws.Range("J22").Value = 440
$ws.Range("K22").Value = 11861.1
$ws.Range("L22").Value = 440
$ws.Range("M22").Value = -11511.1
$ws.Range("N22").Value = -1140
# Row 31
$ws.Range("H31").Value = 3554.1892
$ws.Range("I31").Value = 1965.35
$ws.Range("J31").Value = 5423.4116
$ws.Range("K31").Value = 1965.35
$ws.Range("L31").Value = 5423.4116
$ws.Range("M31").Value = -1670.35
$ws.Range("N31").Value = -6013.4116
# Row 34
$ws.Range("H34").Value = 3554.1892
$ws.Range("I34").Value = 1965.35
$ws.Range("J34").Value = 5423.4116
$ws.Range("K34").Value = 1965.35
$ws.Range("L34").Value = 5423.4116
$ws.Range("M34").Value = -1763.35
$ws.Range("N34").Value = -5827.4116
# Row 97
$ws.Range("H97").Value = 75000
$ws.Range("J97").Value = 75000
$ws.Range("L97").Value = 75000
$ws.Range("N97").Value = -76982
# Row 109
$ws.Range("H109").Value = 108596.8
$ws.Range("J109").Value = 108596.8
$ws.Range("L109").Value = 108596.8
$ws.Range("N109").Value = -110676.8
# Row 138
$ws.Range("H138").Value = 143298
$ws.Range("J138").Value = 143298
$ws.Range("L138").Value = 143298
$ws.Range("N138").Value = -153578
# Row 141
$ws.Range("H141").Value = 238461.47
$ws.Range("J141").Value = 245833.25
$ws.Range("L141").Value = 245833.25
$ws.Range("N141").Value = -256193.25

$ws = $wb.Worksheets.Item("CUL")
# Row 99
$ws.Range("H99").Value = 3286.4285
$ws.Range("I99").Value = 1309.8
$ws.Range("J99").Value = 8228
$ws.Range("K99").Value = 3929.4
$ws.Range("L99").Value = 24684
$ws.Range("M99").Value = -1683.4
$ws.Range("N99").Value = -29176

$ws = $wb.Worksheets.Item("GSM")
# Row 102
$ws.Range("H102").Value = 3967.5918
$ws.Range("I102").Value = 3747.0857
$ws.Range("J102").Value = 4518.857
$ws.Range("K102").Value = 3747.0857
$ws.Range("L102").Value = 4518.857
$ws.Range("M102").Value = -2125.0857
$ws.Range("N102").Value = -7762.857
# Row 107
$ws.Range("H107").Value = 1084.3158
$ws.Range("I107").Value = 1243.2858
$ws.Range("K107").Value = 1243.2858
$ws.Range("M107").Value = 676.7141999999999
# Row 140
$ws.Range("H140").Value = 83569.4
$ws.Range("J140").Value = 83569.4
$ws.Range("L140").Value = 83569.4
$ws.Range("N140").Value = -93929.4

$ws = $wb.Worksheets.Item("LTW")
# Row 68
$ws.Range("H68").Value = 83336310
$ws.Range("I68").Value = 111113780
$ws.Range("J68").Value = 3900
$ws.Range("K68").Value = 111113780
$ws.Range("L68").Value = 3900
$ws.Range("M68").Value = -111113031
$ws.Range("N68").Value = -5398
# Row 71
$ws.Range("H71").Value = 83336310
$ws.Range("I71").Value = 111113780
$ws.Range("J71").Value = 3900
$ws.Range("K71").Value = 555568900
$ws.Range("L71").Value = 19500
$ws.Range("M71").Value = -555565156
$ws.Range("N71").Value = -26988
# Row 105
$ws.Range("H105").Value = 4950
$ws.Range("J105").Value = 4950
$ws.Range("L105").Value = 4950
$ws.Range("N105").Value = -11938
# Row 136
$ws.Range("H136").Value = 9011084
$ws.Range("I136").Value = 12013228
$ws.Range("K136").Value = 36039684
$ws.Range("M136").Value = -36037134

$ws = $wb.Worksheets.Item("WVR")
# Row 107
$ws.Range("H107").Value = 464.78946
$ws.Range("I107").Value = 472.875
$ws.Range("J107").Value = 421.66666
$ws.Range("K107").Value = 1418.625
$ws.Range("L107").Value = 1264.99998
$ws.Range("M107").Value = 501.375
$ws.Range("N107").Value = -5104.999980000001
# Row 122
$ws.Range("H122").Value = 8353
$ws.Range("I122").Value = 5588
$ws.Range("J122").Value = 13729.389
$ws.Range("K122").Value = 16764
$ws.Range("L122").Value = 41188.167
$ws.Range("M122").Value = -14314
$ws.Range("N122").Value = -46088.167
# Row 126
$ws.Range("H126").Value = 3184.1428
$ws.Range("I126").Value = 3058
$ws.Range("J126").Value = 3499.5
$ws.Range("K126").Value = 9174
$ws.Range("L126").Value = 10498.5
$ws.Range("M126").Value = -6704
$ws.Range("N126").Value = -15438.5
# Row 132
$ws.Range("H132").Value = 5953.8945
$ws.Range("I132").Value = 7379
$ws.Range("K132").Value = 22137
$ws.Range("M132").Value = -19607
# Row 136
$ws.Range("H136").Value = 2588.9656
$ws.Range("I136").Value = 2663.6316
$ws.Range("K136").Value = 7990.8948
$ws.Range("M136").Value = -5440.8948
# Row 139
$ws.Range("H139").Value = 68583
$ws.Range("J139").Value = 68583
$ws.Range("L139").Value = 68583
$ws.Range("N139").Value = -78863
